# Issue 47280: LKSM: Trailing whitespace in Source name won't resolve when
# deriving samples.
#
# Reproduces trailing-whitespace values on the sample/source name and
# source-type columns (and turns the previously-boolean E2 cell into the
# literal text "TRUE " with a trailing space) so the regression test data
# matches what LKSM actually saw.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - "SampleSetBVT1" / "a" gain trailing whitespace.
$ws.Range("A2").Value = "SampleSetBVT1  "
$ws.Range("C2").Value = "a "

# Row 3 - "SampleSetBVT2" gains trailing whitespace.
$ws.Range("A3").Value = "SampleSetBVT2   "

# E2 becomes the literal text "TRUE " (trailing space) instead of the
# boolean TRUE it held before.
$ws.Range("E2").Value = "TRUE "

# Leave the selection on E2, matching the authored workbook state.
$ws.Range("E2").Select() | Out-Null
